$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 226.28572
$ws.Range("J9").Value = 286.8
$ws.Range("L9").Value = 286.8
$ws.Range("N9").Value = -624.8
$ws.Range("H100").Value = 886.7143
$ws.Range("I100").Value = 909.5833
$ws.Range("K100").Value = 909.5833
$ws.Range("M100").Value = -368.5833
$ws.Range("H116").Value = 23615186
$ws.Range("I116").Value = 24641720
$ws.Range("K116").Value = 24641720
$ws.Range("M116").Value = -24638278
$ws.Range("H132").Value = 9263.471
$ws.Range("I132").Value = 2749.9167
$ws.Range("J132").Value = 24896
$ws.Range("K132").Value = 8249.750100000001
$ws.Range("L132").Value = 74688
$ws.Range("M132").Value = -5719.750100000001
$ws.Range("N132").Value = -79748
$ws.Range("H137").Value = 13737784
$ws.Range("I137").Value = 1005485.8
$ws.Range("K137").Value = 3016457.4
$ws.Range("M137").Value = -3013907.4
$ws.Range("H138").Value = 6296.6816
$ws.Range("I138").Value = 6295.6665
$ws.Range("J138").Value = 6296.718
$ws.Range("K138").Value = 18886.9995
$ws.Range("L138").Value = 18890.154
$ws.Range("M138").Value = -13746.9995
$ws.Range("N138").Value = -29170.154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14316.323
$ws.Range("I32").Value = 13861
$ws.Range("J32").Value = 18284.143
$ws.Range("K32").Value = 13861
$ws.Range("L32").Value = 18284.143
$ws.Range("M32").Value = -13574
$ws.Range("N32").Value = -18858.143
$ws.Range("H46").Value = 4788776.5
$ws.Range("I46").Value = 19045000
$ws.Range("J46").Value = 36701.832
$ws.Range("K46").Value = 19045000
$ws.Range("L46").Value = 36701.832
$ws.Range("M46").Value = -19044681
$ws.Range("N46").Value = -37339.832
$ws.Range("H61").Value = 6312.1177
$ws.Range("I61").Value = 5824.8184
$ws.Range("J61").Value = 7205.5
$ws.Range("K61").Value = 5824.8184
$ws.Range("L61").Value = 7205.5
$ws.Range("M61").Value = -5612.8184
$ws.Range("N61").Value = -7629.5
$ws.Range("H97").Value = 701
$ws.Range("I97").Value = 375.06668
$ws.Range("J97").Value = 1923.25
$ws.Range("K97").Value = 375.06668
$ws.Range("L97").Value = 1923.25
$ws.Range("M97").Value = 120.93332
$ws.Range("N97").Value = -2915.25
$ws.Range("H132").Value = 13157.389
$ws.Range("I132").Value = 16167.257
$ws.Range("J132").Value = 5331.7334
$ws.Range("K132").Value = 48501.771
$ws.Range("L132").Value = 15995.2002
$ws.Range("M132").Value = -45971.771
$ws.Range("N132").Value = -21055.2002
$ws.Range("H136").Value = 6312.1177
$ws.Range("I136").Value = 5824.8184
$ws.Range("J136").Value = 7205.5
$ws.Range("K136").Value = 17474.4552
$ws.Range("L136").Value = 21616.5
$ws.Range("M136").Value = -14924.4552
$ws.Range("N136").Value = -26716.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3070.0435
$ws.Range("I86").Value = 2055.1177
$ws.Range("K86").Value = 2055.1177
$ws.Range("M86").Value = -932.1176999999998
$ws.Range("H89").Value = 3070.0435
$ws.Range("I89").Value = 2055.1177
$ws.Range("K89").Value = 10275.5885
$ws.Range("M89").Value = -4659.588499999998
$ws.Range("H94").Value = 897.2121
$ws.Range("I94").Value = 806
$ws.Range("K94").Value = 806
$ws.Range("M94").Value = -355
$ws.Range("H99").Value = 2139.7827
$ws.Range("I99").Value = 2236.05
$ws.Range("K99").Value = 2236.05
$ws.Range("M99").Value = -738.0500000000002
$ws.Range("H134").Value = 2864.3572
$ws.Range("I134").Value = 2744.1936
$ws.Range("K134").Value = 8232.5808
$ws.Range("M134").Value = -5697.5808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2714.0833
$ws.Range("I58").Value = 2174.7144
$ws.Range("K58").Value = 2174.7144
$ws.Range("M58").Value = -1971.7144
$ws.Range("H86").Value = 5765.25
$ws.Range("I86").Value = 6208.222
$ws.Range("J86").Value = 5402.8184
$ws.Range("K86").Value = 6208.222
$ws.Range("L86").Value = 5402.8184
$ws.Range("M86").Value = -5085.222
$ws.Range("N86").Value = -7648.8184
$ws.Range("H89").Value = 5765.25
$ws.Range("I89").Value = 6208.222
$ws.Range("J89").Value = 5402.8184
$ws.Range("K89").Value = 31041.11
$ws.Range("L89").Value = 27014.092
$ws.Range("M89").Value = -25425.11
$ws.Range("N89").Value = -38246.092
$ws.Range("H94").Value = 2669.1667
$ws.Range("I94").Value = 325
$ws.Range("J94").Value = 2882.2727
$ws.Range("K94").Value = 325
$ws.Range("L94").Value = 2882.2727
$ws.Range("M94").Value = 126
$ws.Range("N94").Value = -3784.2727
$ws.Range("H97").Value = 50000
$ws.Range("J97").Value = 50000
$ws.Range("L97").Value = 50000
$ws.Range("N97").Value = -51982
$ws.Range("H99").Value = 10878.385
$ws.Range("I99").Value = 11451.583
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 11451.583
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -9953.583000000001
$ws.Range("N99").Value = -6996
$ws.Range("H126").Value = 10878.385
$ws.Range("I126").Value = 11451.583
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 34354.749
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -31884.749
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 40405956
$ws.Range("I132").Value = 53334908
$ws.Range("J132").Value = 2974.25
$ws.Range("K132").Value = 160004724
$ws.Range("L132").Value = 8922.75
$ws.Range("M132").Value = -160002194
$ws.Range("N132").Value = -13982.75
$ws.Range("H136").Value = 2714.0833
$ws.Range("I136").Value = 2174.7144
$ws.Range("K136").Value = 6524.1432
$ws.Range("M136").Value = -3974.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8540.652
$ws.Range("I56").Value = 8540.652
$ws.Range("K56").Value = 8540.652
$ws.Range("M56").Value = -8010.652
$ws.Range("H122").Value = 1988.2667
$ws.Range("I122").Value = 2022
$ws.Range("J122").Value = 1976
$ws.Range("K122").Value = 18198
$ws.Range("L122").Value = 17784
$ws.Range("M122").Value = -15748
$ws.Range("N122").Value = -22684
$ws.Range("H127").Value = 970.8
$ws.Range("J127").Value = 970.8
$ws.Range("L127").Value = 2912.4
$ws.Range("N127").Value = -12832.4
$ws.Range("H131").Value = 10103295
$ws.Range("J131").Value = 12823122
$ws.Range("L131").Value = 38469366
$ws.Range("N131").Value = -38479446

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 401.27585
$ws.Range("I2").Value = 282.14285
$ws.Range("J2").Value = 512.4666999999999
$ws.Range("K2").Value = 282.14285
$ws.Range("L2").Value = 512.4666999999999
$ws.Range("M2").Value = -169.14285
$ws.Range("N2").Value = -738.4666999999999
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -3730
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3064
$ws.Range("N73").ClearContents()
$ws.Range("H113").Value = 5134.8
$ws.Range("I113").Value = 5830.5
$ws.Range("J113").Value = 4671
$ws.Range("K113").Value = 5830.5
$ws.Range("L113").Value = 4671
$ws.Range("M113").Value = -3660.5
$ws.Range("N113").Value = -9011
$ws.Range("H126").Value = 3024.926
$ws.Range("I126").Value = 3065.84
$ws.Range("K126").Value = 9197.52
$ws.Range("M126").Value = -6727.52
$ws.Range("H132").Value = 72212.34
$ws.Range("I132").Value = 97825.95
$ws.Range("J132").Value = 4976.625
$ws.Range("K132").Value = 293477.85
$ws.Range("L132").Value = 14929.875
$ws.Range("M132").Value = -290947.85
$ws.Range("N132").Value = -19989.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1000
$ws.Range("I9").Value = 1000
$ws.Range("K9").Value = 1000
$ws.Range("M9").Value = -776
$ws.Range("H68").Value = 5050.9414
$ws.Range("I68").Value = 4438.0835
$ws.Range("J68").Value = 6521.8
$ws.Range("K68").Value = 4438.0835
$ws.Range("L68").Value = 6521.8
$ws.Range("M68").Value = -3689.0835
$ws.Range("N68").Value = -8019.8
$ws.Range("H71").Value = 5050.9414
$ws.Range("I71").Value = 4438.0835
$ws.Range("J71").Value = 6521.8
$ws.Range("K71").Value = 22190.4175
$ws.Range("L71").Value = 32609
$ws.Range("M71").Value = -18446.4175
$ws.Range("N71").Value = -40097
$ws.Range("H93").Value = 1926.1
$ws.Range("I93").Value = 1681.6666
$ws.Range("K93").Value = 1681.6666
$ws.Range("M93").Value = -433.6666
$ws.Range("H122").Value = 8637.125
$ws.Range("I122").Value = 3699.182
$ws.Range("J122").Value = 12815.385
$ws.Range("K122").Value = 11097.546
$ws.Range("L122").Value = 38446.155
$ws.Range("M122").Value = -8647.545999999998
$ws.Range("N122").Value = -43346.155
$ws.Range("H132").Value = 2575.5
$ws.Range("I132").Value = 2494.7058
$ws.Range("J132").Value = 2747.1875
$ws.Range("K132").Value = 7484.117400000001
$ws.Range("L132").Value = 8241.5625
$ws.Range("M132").Value = -4954.117400000001
$ws.Range("N132").Value = -13301.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5388.391
$ws.Range("I100").Value = 774.38464
$ws.Range("J100").Value = 11386.6
$ws.Range("K100").Value = 1548.76928
$ws.Range("L100").Value = 22773.2
$ws.Range("M100").Value = -1007.76928
$ws.Range("N100").Value = -23855.2
$ws.Range("H107").Value = 2551.3076
$ws.Range("I107").Value = 3324.5715
$ws.Range("K107").Value = 9973.7145
$ws.Range("M107").Value = -8053.7145
$ws.Range("H122").Value = 2244.0476
$ws.Range("I122").Value = 1580.4546
$ws.Range("J122").Value = 2974
$ws.Range("K122").Value = 4741.3638
$ws.Range("L122").Value = 8922
$ws.Range("M122").Value = -2291.3638
$ws.Range("N122").Value = -13822
$ws.Range("H136").Value = 4324.6045
$ws.Range("I136").Value = 2731.375
$ws.Range("J136").Value = 6337.1055
$ws.Range("K136").Value = 8194.125
$ws.Range("L136").Value = 19011.3165
$ws.Range("M136").Value = -5644.125
$ws.Range("N136").Value = -24111.3165
